# Refresh the cryptocurrency price/volume snapshot on the "cryptos" sheet.
# D column values are written with a leading apostrophe so Excel keeps
# them as literal text (matching the source data's formatting, e.g.
# "1.00" / "130.52") instead of auto-converting to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.405.64"
$ws.Range("E2").Value = "  -3.29%  "

$ws.Range("D3").Value = "'2.984.73"
$ws.Range("E3").Value = "  -3.12%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'552.87"
$ws.Range("E5").Value = "  +1.42%  "

$ws.Range("D6").Value = "'130.52"
$ws.Range("E6").Value = "  -6.53%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'2.978.89"
$ws.Range("E8").Value = "  -3.09%  "

$ws.Range("D9").Value = "'0.489"
$ws.Range("E9").Value = "  -2.00%  "

$ws.Range("D10").Value = "'6.00"
$ws.Range("E10").Value = "  -6.58%  "

$ws.Range("D11").Value = "'0.143"
$ws.Range("E11").Value = "  -8.53%  "

$ws.Range("D12").Value = "'0.442"
$ws.Range("E12").Value = "  -3.21%  "

$ws.Range("D13").Value = "'0.0000217"
$ws.Range("E13").Value = "  -3.41%  "

$ws.Range("D14").Value = "'33.78"
$ws.Range("E14").Value = "  -3.59%  "

$ws.Range("D15").Value = "'3.450.06"

$ws.Range("D16").Value = "'61.577.04"
$ws.Range("E16").Value = "  -3.04%  "

$ws.Range("E17").Value = "  -2.95%  "

$ws.Range("D18").Value = "'2.973.36"
$ws.Range("E18").Value = "  -3.52%  "

$ws.Range("D19").Value = "'6.58"
$ws.Range("E19").Value = "  -1.08%  "

$ws.Range("D20").Value = "'478.20"
$ws.Range("E20").Value = "  +0.44%  "

$ws.Range("D21").Value = "'13.06"
$ws.Range("E21").Value = "  -3.12%  "

$ws.Range("D22").Value = "'0.661"
$ws.Range("E22").Value = "  -5.63%  "

$ws.Range("D23").Value = "'6.95"
$ws.Range("E23").Value = "  -2.12%  "

$ws.Range("D24").Value = "'79.50"
$ws.Range("E24").Value = "  +1.02%  "

$ws.Range("D25").Value = "'11.97"
$ws.Range("E25").Value = "  -2.19%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("D27").Value = "'2.70"
$ws.Range("E27").Value = "  -1.05%  "

$ws.Range("D28").Value = "'7.63"
$ws.Range("E28").Value = "  -4.85%  "

$ws.Range("D29").Value = "'0.995"
$ws.Range("E29").Value = "  -0.42%  "

$ws.Range("D30").Value = "'1.89"
$ws.Range("E30").Value = "  -0.42%  "

$ws.Range("D31").Value = "'25.36"
$ws.Range("E31").Value = "  -3.37%  "

$ws.Range("D32").Value = "'1.12"
$ws.Range("E32").Value = "  -2.99%  "

$ws.Range("D33").Value = "'2.29"
$ws.Range("E33").Value = "  -1.01%  "

$ws.Range("D34").Value = "'5.48"
$ws.Range("E34").Value = "  -0.33%  "

$ws.Range("D35").Value = "'54.61"
$ws.Range("E35").Value = "  -6.93%  "

$ws.Range("D36").Value = "'5.84"
$ws.Range("E36").Value = "  -2.99%  "

$ws.Range("D37").Value = "'448.23"
$ws.Range("E37").Value = "  -8.68%  "

$ws.Range("D38").Value = "'3.098.32"
$ws.Range("E38").Value = "  -5.12%  "

$ws.Range("D39").Value = "'0.0787"
$ws.Range("E39").Value = "  -1.52%  "

$ws.Range("D40").Value = "'0.0379"
$ws.Range("E40").Value = "  -6.35%  "

$ws.Range("D41").Value = "'0.116"
$ws.Range("E41").Value = "  -1.97%  "

$ws.Range("D42").Value = "'8.04"
$ws.Range("E42").Value = "  -1.40%  "

$ws.Range("D44").Value = "'2.30"
$ws.Range("E44").Value = "  -11.63%  "

$ws.Range("D45").Value = "'25.33"
$ws.Range("E45").Value = "  -1.62%  "

$ws.Range("D46").Value = "'0.240"
$ws.Range("E46").Value = "  -5.45%  "

$ws.Range("E47").Value = "  -2.03%  "

$ws.Range("D48").Value = "'1.93"
$ws.Range("E48").Value = "  -4.75%  "

$ws.Range("E49").Value = "  +8.84%  "

$ws.Range("D50").Value = "'113.49"
$ws.Range("E50").Value = "  -8.24%  "

$ws.Range("D51").Value = "'0.0₃0478"
$ws.Range("E51").Value = "  -9.91%  "
